$d = $word.ActiveDocument

# Insert a new paragraph right after "Make houses and hotels" containing the
# new to-do item, using Find/Replace with a wildcard paragraph mark (^p) so a
# genuine new <w:p> is created (matching formatting is reapplied below since
# Find/Replace alone does not stamp explicit run properties on the new run).
$newItemText = "If all players say " + [char]0x201C + "no" + [char]0x201D + " to auction, the games stops"
$replacement = "Make houses and hotels^p" + $newItemText

$range = $d.Content
$range.Find.Execute("Make houses and hotels", $true, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)

# Re-apply the same explicit character formatting (Times New Roman, 12pt,
# including complex-script variants) used throughout the rest of the list so
# the newly created run carries explicit rPr like its siblings.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*$newItemText*") {
        $p.Range.Font.Name = "Times New Roman"
        $p.Range.Font.NameBi = "Times New Roman"
        $p.Range.Font.Size = 12
        $p.Range.Font.SizeBi = 12
    }
}
